# continuing refactor of power calculation
#
# Adds a "Link Efficiency (J/bit)" / "Op Efficiency (J/op)" pair of columns to
# both the Detectors and Triggers sheets, moves the old "Compression" column
# on Triggers out of the way (to column I), adds a new "Global" sheet holding
# a single "Year" parameter, and restores the active sheet/selection back to
# Detectors.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New "Global" sheet (after Triggers) holding the projection Year.
# Created first so the "Year" shared string is interned before the
# "Link Efficiency (J/bit)" / "Op Efficiency (J/op)" headers below.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$global = $wb.Worksheets.Add($null, $lastSheet)
$global.Name = "Global"

$global.Range("A1").Value = "Year"
$global.Range("A2").Value = 2028
$global.Range("A3").Select()

# ---------------------------------------------------------------------------
# Detectors (sheet1): add G (Link Efficiency) / H (Op Efficiency) columns
# ---------------------------------------------------------------------------
$detectors = $wb.Worksheets.Item("Detectors")

$detectors.Range("G1").Value = "Link Efficiency (J/bit)"
$detectors.Range("H1").Value = "Op Efficiency (J/op)"

$linkEff = 2.22 * [math]::Pow(10, -11)
$detectors.Range("G2:G20").Value = $linkEff
$detectors.Range("G2:G20").NumberFormat = "0.00E+00"
$detectors.Range("H2:H20").Value = 0

$detectors.Columns.Item(8).ColumnWidth = 19.166666666666668

# ---------------------------------------------------------------------------
# Triggers (sheet2): swap Name/Output columns, insert Link/Op Efficiency
# columns before the old Compression column (now shifted to I)
# ---------------------------------------------------------------------------
$triggers = $wb.Worksheets.Item("Triggers")

for ($r = 1; $r -le 8; $r++) {
    $a = $triggers.Cells.Item($r, 1).Value()
    $b = $triggers.Cells.Item($r, 2).Value()
    $triggers.Cells.Item($r, 1).Value = $b
    $triggers.Cells.Item($r, 2).Value = $a
}

# Preserve the old "Compression" column (G) by moving it to I before
# overwriting G with the new Link Efficiency column.
for ($r = 1; $r -le 8; $r++) {
    $g = $triggers.Cells.Item($r, 7).Value()
    $triggers.Cells.Item($r, 9).Value = $g
}

$triggers.Range("G1").Value = "Link Efficiency (J/bit)"
$triggers.Range("H1").Value = "Op Efficiency (J/op)"
$triggers.Range("G1:H1").Font.Color = 0

$trigEff = 2.5 * [math]::Pow(10, -11)
$triggers.Range("G2:G8").Value = $trigEff
$triggers.Range("G2:G8").NumberFormat = "0.00E+00"
$triggers.Range("H2:H8").Value = 0

$triggers.Range("H6").Formula = "=120000/(40000000)"
$triggers.Range("H7").Formula = "=1600000/(100000)"

$triggers.Columns.Item(7).ColumnWidth = 20.830729166666668
$triggers.Columns.Item(8).ColumnWidth = 21.166666666666668

# ---------------------------------------------------------------------------
# Restore selections / active sheet
# ---------------------------------------------------------------------------
$triggers.Range("D15").Select()
$detectors.Range("J19").Select()
$detectors.Activate()
